# "preparing for final checkin"
# - Replace the placeholder name (Rahul Ranjan) on Sheet2 with the final
#   test-account name (Test User); postcode is left untouched.
# - Leave Sheet1 as the active/selected tab (was Sheet2) and update the
#   last selected cell remembered on each sheet.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")

# Sheet2 data: firstname/lastname columns (C2/D2). Postcode (E2) stays CF117JA.
$ws2.Range("C2").Value = "Test"
$ws2.Range("D2").Value = "User"

# Remember the last selection on Sheet2 before leaving it.
$ws2.Range("C7").Select()

# Make Sheet1 the active sheet/tab and set its remembered selection.
$ws1.Activate()
$ws1.Range("A15").Select()
